$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44767
$ws.Range("J2").Value = 50
$ws.Range("D3").Value = 44756
$ws.Range("J3").Value = 80
$ws.Range("D4").Value = 44845
$ws.Range("J4").Value = 20
$ws.Range("K4").Value = 16000
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 16000
$ws.Range("P4").Value = 1067
$ws.Range("D5").Value = 44839
$ws.Range("J5").Value = 80
$ws.Range("K5").Value = 16000
$ws.Range("L5").Value = 16000
$ws.Range("M5").Value = 16000
$ws.Range("P5").Value = 1067
$ws.Range("D6").Value = 44837
$ws.Range("J6").Value = 80
$ws.Range("K6").Value = 16000
$ws.Range("L6").Value = 16000
$ws.Range("M6").Value = 16000
$ws.Range("P6").Value = 1067
$ws.Range("D7").Value = 44830
$ws.Range("J7").Value = 25
$ws.Range("K7").Value = 12000
$ws.Range("L7").Value = 12000
$ws.Range("M7").Value = 12000
$ws.Range("P7").Value = 800
$ws.Range("D8").Value = 44819
$ws.Range("J8").Value = 100
$ws.Range("D9").Value = 44525
$ws.Range("J9").Value = 40
$ws.Range("K9").Value = 8000
$ws.Range("L9").Value = 8000
$ws.Range("M9").Value = 8000
$ws.Range("P9").Value = 533
$ws.Range("D10").Value = 44776
$ws.Range("J10").Value = 80
$ws.Range("D11").Value = 44827
$ws.Range("D12").Value = 44824
$ws.Range("J12").Value = 20
$ws.Range("K12").Value = 20000
$ws.Range("L12").Value = 20000
$ws.Range("M12").Value = 20000
$ws.Range("P12").Value = 1333
$ws.Range("D13").Value = 44755
$ws.Range("J13").Value = 50
$ws.Range("D14").Value = 44508
$ws.Range("J14").Value = 40
$ws.Range("K14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 10000
$ws.Range("P14").Value = 667
$ws.Range("D15").Value = 45134
$ws.Range("J15").Value = 5
$ws.Range("D16").Value = 44841
$ws.Range("J16").Value = 20
$ws.Range("K16").Value = 16000
$ws.Range("L16").Value = 16000
$ws.Range("M16").Value = 16000
$ws.Range("P16").Value = 1067
$ws.Range("D17").Value = 44813
$ws.Range("J17").Value = 20
$ws.Range("K17").Value = 20000
$ws.Range("L17").Value = 20000
$ws.Range("M17").Value = 20000
$ws.Range("P17").Value = 1333
$ws.Range("D18").Value = 44518
$ws.Range("J18").Value = 50
$ws.Range("K18").Value = 10000
$ws.Range("L18").Value = 10000
$ws.Range("M18").Value = 10000
$ws.Range("P18").Value = 667
$ws.Range("D19").Value = 44812
$ws.Range("J19").Value = 80
$ws.Range("K19").Value = 20000
$ws.Range("L19").Value = 20000
$ws.Range("M19").Value = 20000
$ws.Range("P19").Value = 1333
$ws.Range("D20").Value = 44769
$ws.Range("J20").Value = 50
$ws.Range("D21").Value = 44825
$ws.Range("J21").Value = 30
$ws.Range("D22").Value = 44838
$ws.Range("J22").Value = 10
$ws.Range("K22").Value = 20000
$ws.Range("L22").Value = 20000
$ws.Range("M22").Value = 20000
$ws.Range("P22").Value = 1333
$ws.Range("D23").Value = 44757
$ws.Range("D24").Value = 44826
$ws.Range("J24").Value = 50
$ws.Range("D25").Value = 44771
$ws.Range("J25").Value = 40
$ws.Range("D26").Value = 44811
$ws.Range("J26").Value = 30
$ws.Range("K26").Value = 20000
$ws.Range("L26").Value = 20000
$ws.Range("M26").Value = 20000
$ws.Range("P26").Value = 1333
$ws.Range("D27").Value = 44749
$ws.Range("J27").Value = 50